# ------------------------------------------------------------------
# Add 2022-Q3 data:
#  - "总计" (summary) sheet: shift quarterly totals down one slot and
#    insert the new 2022-Q3 figures at the top (2020-Q4 falls into a
#    newly appended row).
#  - Insert a brand-new worksheet "2022-Q3" (fund holdings detail)
#    positioned right after "总计" and before "2022-Q2".
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "总计" (summary) sheet
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Quarterly figures, newest first, matching the final layout of rows 2-9.
$quarters = @(
    @("2022-Q3", 18, 4.95),
    @("2022-Q2", 57, 12.68),
    @("2022-Q1", 20, 5.59),
    @("2021-Q4", 33, 9.18),
    @("2021-Q3", 27, 6.59),
    @("2021-Q2", 18, 5.18),
    @("2021-Q1", 3, 0.33),
    @("2020-Q4", 2, 0.09)
)

# Append a new row 9 (copy formatting from row 8) to hold the quarter
# that falls off the bottom of the table (2020-Q4).
$total.Range("A8:D8").Copy()
$total.Range("A9:D9").PasteSpecial(-4122)
$total.Cells.Item(9, 1).Value = 7

# Rewrite the date/count/value columns for rows 2-9 (index column A is
# left untouched - it already holds the correct sequential values).
for ($i = 0; $i -lt $quarters.Length; $i++) {
    $r = $i + 2
    $q = $quarters[$i]
    $total.Cells.Item($r, 2).Value = $q[0]
    $total.Cells.Item($r, 3).Value = $q[1]
    $total.Cells.Item($r, 4).Value = $q[2]
}

# ---------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet, positioned before "2022-Q2"
# ---------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (row 1), columns B..H
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3Sheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Fund holding detail rows: index, code, name, scale, position%, ratio%, marketValue, rank
$q3Data = @(
    @(0, "001985", "富国低碳新经济混合A", "22.20", "92.54", "6.33", "1.4053", 3),
    @(1, "519035", "富国天博创新混合", "18.16", "91.30", "5.70", "1.0351", 2),
    @(2, "011357", "华泰柏瑞品质成长混合A", "21.25", "74.12", "2.85", "0.6056", 7),
    @(3, "000513", "富国高端制造行业股票A", "7.01", "91.41", "5.88", "0.4122", 3),
    @(4, "006218", "富国生物医药科技混合A", "7.10", "89.32", "4.27", "0.3032", 7),
    @(5, "008138", "富国龙头优势混合", "4.41", "92.58", "5.95", "0.2624", 2),
    @(6, "009990", "华泰柏瑞品质优选混合A", "9.06", "68.27", "2.75", "0.2492", 6),
    @(7, "011921", "富国均衡成长三年持有期混合A", "6.14", "90.69", "3.66", "0.2247", 5),
    @(8, "010122", "华泰柏瑞优势领航混合A", "3.95", "79.94", "3.77", "0.1489", 2),
    @(9, "100016", "富国天源沪港深平衡混合A", "4.99", "70.11", "1.81", "0.0903", 7),
    @(10, "009991", "华泰柏瑞品质优选混合C", "2.49", "68.27", "2.75", "0.0685", 6),
    @(11, "011308", "富国生物医药科技混合C", "1.52", "89.32", "4.27", "0.0649", 7),
    @(12, "011358", "华泰柏瑞品质成长混合C", "1.10", "74.12", "2.85", "0.0314", 7),
    @(13, "011922", "富国均衡成长三年持有期混合C", "0.44", "90.69", "3.66", "0.0161", 5),
    @(14, "011306", "富国低碳新经济混合C", "0.25", "92.54", "6.33", "0.0158", 3),
    @(15, "010123", "华泰柏瑞优势领航混合C", "0.39", "79.94", "3.77", "0.0147", 2),
    @(16, "014930", "富国高端制造行业股票C", "0.01", "91.41", "5.88", "0.0006", 3),
    @(17, "014931", "富国天源沪港深平衡混合C", "0.00", "70.11", "1.81", "0", 7)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = $q3Data[$i]
    $r = $i + 2

    $q3Sheet.Cells.Item($r, 1).Value = $row[0]
    $q3Sheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3Sheet.Cells.Item($r, 3).Value = $row[2]
    $q3Sheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3Sheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3Sheet.Cells.Item($r, 6).Value = "'" + $row[5]

    if ($i -eq ($q3Data.Length - 1)) {
        # last row's market-value is exactly 0, stored as a real number
        $q3Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q3Sheet.Cells.Item($r, 7).Value = "'" + $row[6]
    }

    $q3Sheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------
# 3) Apply formatting to match the other quarterly sheets:
#    header row (B1:H1) and index column (A2:A19) use the bold /
#    bordered style already present elsewhere in the workbook.
# ---------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("2022-Q2")

$styleSource.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)

$lastRow = $q3Data.Length + 1
$styleSource.Range("A2").Copy()
$q3Sheet.Range("A2:A" + $lastRow).PasteSpecial(-4122)

$q3Sheet.Range("A1").Select()
